$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 104 - this shifts existing rows 104:126 down to 105:127
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with a new weekly price record for
# "Rabanito" (same market/region/category as its neighbours), dated 2023-09-04
# (serial 45173) with a Volumen of 50.
$ws.Range("A104").Value = 10
$ws.Range("B104").Value = "Vega Modelo de Temuco"
$ws.Range("C104").Value = "La Araucanía"
$ws.Range("D104").NumberFormat = $ws.Range("D105").NumberFormat
$ws.Range("D104").Value = 45173
$ws.Range("E104").Value = 9
$ws.Range("F104").Value = 300000001
$ws.Range("G104").Value = "Rabanito"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 50
$ws.Range("K104").Value = 8000
$ws.Range("L104").Value = 8000
$ws.Range("M104").Value = 8000
$ws.Range("N104").Value = "`$/docena de paquetes"
$ws.Range("O104").Value = "Provincia de Cautín"
$ws.Range("P104").Value = 667
$ws.Range("Q104").Value = 12
$ws.Range("R104").Value = "Hortaliza"
